$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.067.11'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '3.103.72'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.098.52'
$ws.Range("E8").Value = '  +2.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("E10").Value = '  -4.14%  '
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.22'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '3.618.51'
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("D17").Value = '67.086.49'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").Value = '3.105.39'
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '485.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.93%  '
$ws.Range("E22").Value = '  +1.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.30%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.55%  '
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.74%  '
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.114'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.992'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.68'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.316'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("D45").Value = '2.849.70'
$ws.Range("E45").Value = '  +4.50%  '
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '384.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.04%  '
